# Commit: "change opcode of some instructions"
#
# The HALFMIPS workbook documents the MIPS-I and HALF-MIPS opcode maps.
# This edit adds the two missing COP0 memory instructions (LWC0 / SWC0)
# to both opcode tables, and re-shuffles a handful of cells on the
# HALF-MIPS sheet so the "half" instruction variants line up with the
# now-present LWC0/SWC0 slots.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "MIPS-I": opcode map — bits 31..29 = "110"(row10)/"111"(row11)
# The COP0 row previously had "*" placeholders in the LW/SW columns
# (bits 28..26 = "000"); fill in the real opcode names.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("MIPS-I")

$ws1.Range("D10").Value = "LWC0"
$ws1.Range("D11").Value = "SWC0"

# Match the italic styling already used for the SPECIAL/REGIMM headers
# on the HALF-MIPS sheet.
$ws1.Range("D4:E4").Font.Italic = $true

# ---------------------------------------------------------------------
# Sheet "HALF-MIPS": opcode map mirrors MIPS-I but with half-word/byte
# "H" instructions filling most of the COP0 row. Adding LWC0/SWC0 bumps
# HLB/HLBU one slot to the right within that row.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("HALF-MIPS")

$ws2.Range("D10").Value = "LWC0"
$ws2.Range("H10").Value = "HLB"
$ws2.Range("J10").Value = "HLBU"
$ws2.Range("D11").Value = "SWC0"

# HSLT / HSLTU move from the SPECIAL2 "111" row (24) up into the "101"
# row (22), where they replace two "*" placeholders.
$ws2.Range("J22").Value = "HSLT"
$ws2.Range("K22").Value = "HSLTU"

# Row 24 (SPECIAL2, bits 28..26 = "111") shifts left by one to close the
# gap left by HSLT/HSLTU, leaving the last two slots as "*".
$ws2.Range("E24").Value = "*"
$ws2.Range("F24").Value = "HSRL"
$ws2.Range("G24").Value = "HSRA"
$ws2.Range("H24").Value = "HSLLV"
$ws2.Range("I24").Value = "*"

# ---------------------------------------------------------------------
# Sheet "Sheet3": give it an explicit page setup (paper size / orientation)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1
